$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.267.92'
$ws.Range("E2").Value = '  -5.95%  '
$ws.Range("D3").Value = '2.456.67'
$ws.Range("E3").Value = '  -8.43%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '539.47'
$ws.Range("E5").Value = '  -2.82%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.77'
$ws.Range("E6").Value = '  -6.92%  '
$ws.Range("E7").Value = '  -0.18%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.568'
$ws.Range("E8").Value = '  -4.26%  '
$ws.Range("D9").Value = '2.473.66'
$ws.Range("E9").Value = '  -7.93%  '
$ws.Range("E10").Value = '  -6.05%  '
$ws.Range("E11").Value = '  -2.50%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.29'
$ws.Range("E12").Value = '  -1.55%  '
$ws.Range("E13").Value = '  -4.25%  '
$ws.Range("D14").Value = '2.897.37'
$ws.Range("E14").Value = '  -8.22%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '24.10'
$ws.Range("E15").Value = '  -8.51%  '
$ws.Range("D16").Value = '59.221.88'
$ws.Range("E16").Value = '  -5.83%  '
$ws.Range("E17").Value = '  -6.14%  '
$ws.Range("D18").Value = '2.525.15'
$ws.Range("E18").Value = '  -5.87%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.15'
$ws.Range("E19").Value = '  -6.46%  '
$ws.Range("E20").Value = '  -6.13%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '324.41'
$ws.Range("E21").Value = '  -6.11%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.73'
$ws.Range("E23").Value = '  -9.13%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.461'
$ws.Range("E24").Value = '  -9.42%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '60.57'
$ws.Range("E25").Value = '  -4.41%  '
$ws.Range("E26").Value = '  -3.96%  '
$ws.Range("E27").Value = '  -2.01%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.70'
$ws.Range("E28").Value = '  -6.35%  '
$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.82'
$ws.Range("E29").Value = '  -6.58%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.74'
$ws.Range("E30").Value = '  -7.28%  '
$ws.Range("B31").Value = 'Fetch.AI'
$ws.Range("C31").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.26'
$ws.Range("E31").Value = '  -7.78%  '
$ws.Range("D32").Value = '0.0₃0774'
$ws.Range("E32").Value = '  -10.01%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.998'
$ws.Range("E33").Value = '  -0.10%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '158.45'
$ws.Range("E34").Value = '  -4.17%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.56'
$ws.Range("E35").Value = '  -6.80%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.37'
$ws.Range("E36").Value = '  -7.63%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '18.39'
$ws.Range("E37").Value = '  -5.92%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.75'
$ws.Range("E38").Value = '  -1.85%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.92'
$ws.Range("E39").Value = '  -6.84%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '320.02'
$ws.Range("E40").Value = '  -8.23%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '36.72'
$ws.Range("E41").Value = '  -4.10%  '
$ws.Range("E42").Value = '  -12.76%  '
$ws.Range("E43").Value = '  -7.74%  '
$ws.Range("E44").Value = '  -0.32%  '
$ws.Range("E45").Value = '  -2.76%  '
$ws.Range("E46").Value = '  -5.26%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0940'
$ws.Range("E47").Value = '  -3.50%  '
$ws.Range("E48").Value = '  -6.68%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '19.03'
$ws.Range("E49").Value = '  -8.82%  '
$ws.Range("B50").Value = 'VeChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0229'
$ws.Range("E50").Value = '  -5.22%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '18.56'
$ws.Range("E51").Value = '  -9.15%  '
